# davco 65+ initial commit
# Add the new "FiftyForward Madison Station Senior Center" record (row 19)
# and unhide / resize the supporting data columns (D, E, F, I, J, K, L) that
# back the visible summary columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unhide helper columns and restore their working widths ---------------
$ws.Columns.Item(4).Hidden  = $false
$ws.Columns.Item(5).Hidden  = $false
$ws.Columns.Item(5).ColumnWidth = 37.85
$ws.Columns.Item(6).Hidden  = $false
$ws.Columns.Item(9).Hidden  = $false
$ws.Columns.Item(10).Hidden = $false
$ws.Columns.Item(11).Hidden = $false
$ws.Columns.Item(12).Hidden = $false

# --- New row of data: FiftyForward Madison Station Senior Center ----------
$row = 19
$ws.Range("A$row").Value = "69-72"
$ws.Range("B$row").Value = "FiftyForward Madison Station Senior Center"
$ws.Range("C$row").Value = "Davidson"
$ws.Range("D$row").Value = "530 Madison Station Blvd, Madison, TN 37115"
$ws.Range("F$row").Value = "Fifty Forward Madison Station 5 Mile Radius"
$ws.Range("E$row").Value = "FiftyForward Madison Station"
$ws.Range("G$row").Value = "2020 ACS"
$ws.Range("H$row").Value = "2020 ACS"
$ws.Range("I$row").Value = "Done w links"
$ws.Range("J$row").Value = "Done w links"
$ws.Range("K$row").Value = "Done w links"
$ws.Range("L$row").Value = "Done w links"
$ws.Range("M$row").Value = "5m map"
$ws.Range("N$row").Value = "https://reports.mysidewalk.com/b797651eb3"

# Hyperlink the sharing-link cell, matching the style used by the rows above.
$ws.Hyperlinks.Add($ws.Range("N$row"), "https://reports.mysidewalk.com/b797651eb3")
$ws.Range("N$row").Style = "Hyperlink"

# --- Selection / scroll position left by the editor ------------------------
$ws.Range("M25").Select()
